$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2025-10-04 18:55:06", "Noah", 8450689526, "Hey man what’s up?"),
    @("2025-10-04 18:56:34", "Noah", 8450689526, "Hey man what’s up?"),
    @("2025-10-04 19:25:39", "Noah", 8450689526, "Hey man what’s up?"),
    @("2025-10-04 19:27:47", "Noah", 8450689526, "Hey man what’s up?"),
    @("2025-10-04 19:29:36", "Noah", 8450689526, "Test message"),
    @("2025-10-04 19:30:42", "Noah", 8450689526, "Test message"),
    @("2025-10-04 19:31:22", "Noah", 8450689526, "Test message")
)

$startRow = 57
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    # Column D ("Phone") is stored as text in this sheet (e.g. D6) -
    # copy from an existing text-typed phone cell so the type/format matches,
    # then the value already equals the target phone number.
    $ws.Range("D6").Copy($ws.Cells.Item($r, 4))
    $ws.Cells.Item($r, 5).Value = $data[3]
    # Columns F ("Media") and G ("Channel") are blank text cells for these rows
    # (inline-string cells with no text) - copy from existing blank text cells
    # (F2/G2) rather than assigning "" (which would clear the cell entirely).
    $ws.Range("F2").Copy($ws.Cells.Item($r, 6))
    $ws.Range("G2").Copy($ws.Cells.Item($r, 7))
}
